# "cambiati part number led" - change the LED part number and add a
# "spesa" (cost) helper column that multiplies quantity by the unit
# price of the LED (185), replacing the old hard-coded *180 helper.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The LED part was re-sourced; drop the old Mouser hyperlink (and its
# "859-" distributor prefix) and replace the part number text with the
# bare manufacturer part number.
$ws.Hyperlinks.Delete()
$ws.Range("D5").Value = "LTST-C230TBKT"

# New unit price reference cell for the LED (185), and a "spesa" column
# (F) that multiplies each line's quantity by that unit price.
$ws.Range("H1").Value = 185
$ws.Range("F2").Formula = '=A2*$H$1'
$ws.Range("F3:F7").Formula = '=A3*$H$1'

# Widen column F to fit the new values and move the selection where the
# author last clicked.
$ws.Columns.Item(6).ColumnWidth = 14
$ws.Range("G5").Select()
